# Refactor export service to support chunk query with only one iteration
#
# The export template's second row used to hold double-templated
# placeholders ("{{_.EmployeeNo}}", "{{_.Salary}}", "{{_.FromDate}}",
# "{{_.ToDate}}", "{{_.create_at}}") mirroring the header row's field
# names. Now that the exporter streams the query in a single chunked
# iteration, that row instead carries the literal binding keys used to
# look the values up (plain field names, with "create_at" replacing the
# old "{{_.create_at}}" placeholder), and the header row gets a green
# highlight fill so it's visually distinguishable as the header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: replace the old "{{_.Field}}" template placeholders with the
# plain field/key names (same names as the header row, except the last
# column which becomes "create_at").
$ws.Range("A2").Value = "EmployeeNo"
$ws.Range("B2").Value = "Salary"
$ws.Range("C2").Value = "FromDate"
$ws.Range("D2").Value = "ToDate"
$ws.Range("E2").Value = "create_at"

# Row 1 (header): highlight with a green fill, keep centered + bordered.
$ws.Range("A1:E1").Interior.Color = 5287936
$ws.Range("A1:E1").HorizontalAlignment = -4108
$ws.Range("A1:E1").VerticalAlignment = -4108

# Restore the cursor/selection to where the author last left it.
$ws.Range("E7").Select() | Out-Null
